$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'لا حزن يبقى والقرآن الكريم يتلى الشيخ ياسر الدوسري.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'ধুমপান করে নামায পরলে নামায হবে কি-.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Subtitle) - Muhammad Al Muqit - REMINDER.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'শুধু আল্লাহকে দেখানোর জন্য কোনো কাজ করেছেন কখনো- মিশারী আল আফাসীর তাকওয়া জাগানিয়া সংগীত.mp4',
    'Facebook_4.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '১০ শ্রেণীর লোক মারা গেলে শাহাদাতের মর্যাদা পাবে / আলোচক ড আবু বকর মুহাম্মাদ যাকারিয়া হাফিযাহুল্লাহ.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'لا حزن يبقى والقرآن الكريم يتلى الشيخ ياسر الدوسري.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'ধুমপান করে নামায পরলে নামায হবে কি-.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Subtitle) - Muhammad Al Muqit - REMINDER.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'শুধু আল্লাহকে দেখানোর জন্য কোনো কাজ করেছেন কখনো- মিশারী আল আফাসীর তাকওয়া জাগানিয়া সংগীত.mp4',
    'Facebook_4.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '১০ শ্রেণীর লোক মারা গেলে শাহাদাতের মর্যাদা পাবে / আলোচক ড আবু বকর মুহাম্মাদ যাকারিয়া হাফিযাহুল্লাহ.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'لا حزن يبقى والقرآن الكريم يتلى الشيخ ياسر الدوسري.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'ধুমপান করে নামায পরলে নামায হবে কি-.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Subtitle) - Muhammad Al Muqit - REMINDER.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'শুধু আল্লাহকে দেখানোর জন্য কোনো কাজ করেছেন কখনো- মিশারী আল আফাসীর তাকওয়া জাগানিয়া সংগীত.mp4',
    'Facebook_4.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '১০ শ্রেণীর লোক মারা গেলে শাহাদাতের মর্যাদা পাবে / আলোচক ড আবু বকর মুহাম্মাদ যাকারিয়া হাফিযাহুল্লাহ.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'لا حزن يبقى والقرآن الكريم يتلى الشيخ ياسر الدوسري.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'ধুমপান করে নামায পরলে নামায হবে কি-.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Subtitle) - Muhammad Al Muqit - REMINDER.mp4',
    '#مشاري_راشد_العفاسي حب اليتامى - Alafasy Hob Alyatamah - YouTube.MP4',
    'শুধু আল্লাহকে দেখানোর জন্য কোনো কাজ করেছেন কখনো- মিশারী আল আফাসীর তাকওয়া জাগানিয়া সংগীত.mp4',
    'Facebook_4.mp4',
    '১৭ শ্রেণীর লোক মুসলমান হয়েও জান্নাতে যেতে পারবে না! -শায়খ আহমাদুল্লাহ.mp4',
    '১০ শ্রেণীর লোক মারা গেলে শাহাদাতের মর্যাদা পাবে / আলোচক ড আবু বকর মুহাম্মাদ যাকারিয়া হাফিযাহুল্লাহ.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'No sadness remains and the Noble Qur’an is recited by Sheikh Yasser Al-Dosari.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'If you pray after smoking, what will be the prayers -. Mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Thipttel) - Muhammad Al Muqit - Reminder.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Has he ever done anything just to show Allah - Mishari Al Afasi''s Taqwa Jagania Sangeet.mp4',
    'Facebook_4.mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    'If class 10 people die, they will get the status of martyrdom / Dr. Abu Bakr Muhammad Zakaria Hafizahullah.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'No sadness remains and the Noble Qur’an is recited by Sheikh Yasser Al-Dosari.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'If you pray after smoking, what will be the prayers -. Mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Thipttel) - Muhammad Al Muqit - Reminder.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Has he ever done anything just to show Allah - Mishari Al Afasi''s Taqwa Jagania Sangeet.mp4',
    'Facebook_4.mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    'If class 10 people die, they will get the status of martyrdom / Dr. Abu Bakr Muhammad Zakaria Hafizahullah.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'No sadness remains and the Noble Qur’an is recited by Sheikh Yasser Al-Dosari.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'If you pray after smoking, what will be the prayers -. Mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Thipttel) - Muhammad Al Muqit - Reminder.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Has he ever done anything just to show Allah - Mishari Al Afasi''s Taqwa Jagania Sangeet.mp4',
    'Facebook_4.mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    'If class 10 people die, they will get the status of martyrdom / Dr. Abu Bakr Muhammad Zakaria Hafizahullah.mp4',
    'Emotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur RehmanEmotional Nasheed - Har Ek Janib Udasi Hai - Na Jane Kya Ye Alum Hay- Atiq Ur Rehman.mp4',
    '❤ Allah tar rasul ke koto valobashten! 😞 Very Emotional Waz with English and bangla subtitle',
    'Amazing Drone Footage Of Al Masjid Al Haram - YouTube.mp4',
    'No sadness remains and the Noble Qur’an is recited by Sheikh Yasser Al-Dosari.mp4',
    'Jore Amin Bola By Mufti Kazi Ibrahim - YouTube.MP4',
    'If you pray after smoking, what will be the prayers -. Mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Bangla Nasheed - Bab Al Raja (Bangla Thipttel) - Muhammad Al Muqit - Reminder.mp4',
    '#Mishary_Rashed_Al-Afasy love for orphans - Alafasy Hob Alyatamah - YouTube.MP4',
    'Has he ever done anything just to show Allah - Mishari Al Afasi''s Taqwa Jagania Sangeet.mp4',
    'Facebook_4.mp4',
    '16 class people will not be able to go to paradise even if they are Muslims! -Shaykh Ahmadullah.mp4',
    'If class 10 people die, they will get the status of martyrdom / Dr. Abu Bakr Muhammad Zakaria Hafizahullah.mp4'
)

$startRow = 87
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}
